$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns D (Description), E (Reference), F (Owner)
$ws.Range("D1").Value = "Description"
$ws.Range("E1").Value = "Reference"
$ws.Range("F1").Value = "Owner"

$ws.Range("D2").Value = "Introduction to concepts of biology, like gene, DNA etc"
$ws.Range("E2").Value = "SVN reporistory under docs/Induction"
$ws.Range("F2").Value = "Madhumita Shrikhande"

$ws.Range("D3").Value = "Concepts about caBIG, their objective, what is caGRID, how the data is stored etc"
$ws.Range("F3").Value = "Chandrakant Talele"

$ws.Range("D4").Value = "Demonstration of all the Admin module features"
$ws.Range("E4").Value = "User Manual in SVN repository"
$ws.Range("F4").Value = "Pooja Arora"

$ws.Range("D5").Value = "Demonstration of all the Desktop Application features"
$ws.Range("E5").Value = "User Manual in SVN repository"
$ws.Range("F5").Value = "Pooja Arora"

$ws.Range("D6").Value = "Demonstration of all the web application features"
$ws.Range("E6").Value = "User Manual in SVN repository"
$ws.Range("F6").Value = "Pooja Arora"

$ws.Range("D7").Value = "Explanation of the design and overall architechture of the entire application"
$ws.Range("E7").Value = "Design and Archtechture Document in SVN repository"
$ws.Range("F7").Value = "Chandrakant Talele"

$ws.Range("E8").Value = "Design and Archtechture Document in SVN repository"
$ws.Range("F8").Value = "Chandrakant Talele"

$ws.Range("F9").Value = "Chandrakant Talele"

# Column widths for the new D and E columns
$ws.Columns.Item(4).ColumnWidth = 48.8
$ws.Columns.Item(5).ColumnWidth = 16.67

# Update the cursor/selection to D8, matching the edited sheet's last view state
$ws.Range("D8").Select() | Out-Null
